$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2021-Q1")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)
$newSheet.Name = "2022-Q1"

# Row 19 doesnt exist in source (only 18 rows); seed its A-column style from A18
$newSheet.Range("A18").Copy()
$newSheet.Range("A19").PasteSpecial(-4122)

$newSheet.Cells.Item(1,2).Value2 = "基金代码"
$newSheet.Cells.Item(1,3).Value2 = "基金名称"
$newSheet.Cells.Item(1,4).Value2 = "基金规模"
$newSheet.Cells.Item(1,5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1,6).Value2 = "仓位占比"
$newSheet.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value2 = "仓位排名"

# row 2
$newSheet.Cells.Item(2,1).Value2 = 0
$newSheet.Cells.Item(2,2).NumberFormat = "@"
$newSheet.Cells.Item(2,2).Value2 = "501092"
$newSheet.Cells.Item(2,2).ClearFormats()
$newSheet.Cells.Item(2,3).NumberFormat = "@"
$newSheet.Cells.Item(2,3).Value2 = "交银施罗德瑞思三年封闭运作混合"
$newSheet.Cells.Item(2,3).ClearFormats()
$newSheet.Cells.Item(2,4).NumberFormat = "@"
$newSheet.Cells.Item(2,4).Value2 = "64.69"
$newSheet.Cells.Item(2,4).ClearFormats()
$newSheet.Cells.Item(2,5).NumberFormat = "@"
$newSheet.Cells.Item(2,5).Value2 = "88.90"
$newSheet.Cells.Item(2,5).ClearFormats()
$newSheet.Cells.Item(2,6).NumberFormat = "@"
$newSheet.Cells.Item(2,6).Value2 = "2.15"
$newSheet.Cells.Item(2,6).ClearFormats()
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value2 = "1.3908"
$newSheet.Cells.Item(2,7).ClearFormats()
$newSheet.Cells.Item(2,8).Value2 = 6

# row 3
$newSheet.Cells.Item(3,1).Value2 = 1
$newSheet.Cells.Item(3,2).NumberFormat = "@"
$newSheet.Cells.Item(3,2).Value2 = "010671"
$newSheet.Cells.Item(3,2).ClearFormats()
$newSheet.Cells.Item(3,3).NumberFormat = "@"
$newSheet.Cells.Item(3,3).Value2 = "景顺长城大中华混合(QDII)美元"
$newSheet.Cells.Item(3,3).ClearFormats()
$newSheet.Cells.Item(3,4).NumberFormat = "@"
$newSheet.Cells.Item(3,4).Value2 = "10.35"
$newSheet.Cells.Item(3,4).ClearFormats()
$newSheet.Cells.Item(3,5).NumberFormat = "@"
$newSheet.Cells.Item(3,5).Value2 = "82.59"
$newSheet.Cells.Item(3,5).ClearFormats()
$newSheet.Cells.Item(3,6).NumberFormat = "@"
$newSheet.Cells.Item(3,6).Value2 = "4.97"
$newSheet.Cells.Item(3,6).ClearFormats()
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value2 = "0.5144"
$newSheet.Cells.Item(3,7).ClearFormats()
$newSheet.Cells.Item(3,8).Value2 = 6

# row 4
$newSheet.Cells.Item(4,1).Value2 = 2
$newSheet.Cells.Item(4,2).NumberFormat = "@"
$newSheet.Cells.Item(4,2).Value2 = "262001"
$newSheet.Cells.Item(4,2).ClearFormats()
$newSheet.Cells.Item(4,3).NumberFormat = "@"
$newSheet.Cells.Item(4,3).Value2 = "景顺长城大中华混合(QDII)"
$newSheet.Cells.Item(4,3).ClearFormats()
$newSheet.Cells.Item(4,4).NumberFormat = "@"
$newSheet.Cells.Item(4,4).Value2 = "10.35"
$newSheet.Cells.Item(4,4).ClearFormats()
$newSheet.Cells.Item(4,5).NumberFormat = "@"
$newSheet.Cells.Item(4,5).Value2 = "82.59"
$newSheet.Cells.Item(4,5).ClearFormats()
$newSheet.Cells.Item(4,6).NumberFormat = "@"
$newSheet.Cells.Item(4,6).Value2 = "4.97"
$newSheet.Cells.Item(4,6).ClearFormats()
$newSheet.Cells.Item(4,7).NumberFormat = "@"
$newSheet.Cells.Item(4,7).Value2 = "0.5144"
$newSheet.Cells.Item(4,7).ClearFormats()
$newSheet.Cells.Item(4,8).Value2 = 6

# row 5
$newSheet.Cells.Item(5,1).Value2 = 3
$newSheet.Cells.Item(5,2).NumberFormat = "@"
$newSheet.Cells.Item(5,2).Value2 = "862001"
$newSheet.Cells.Item(5,2).ClearFormats()
$newSheet.Cells.Item(5,3).NumberFormat = "@"
$newSheet.Cells.Item(5,3).Value2 = "光大阳光香港精选混合型集合资产管理计划（QDII）A 人民币"
$newSheet.Cells.Item(5,3).ClearFormats()
$newSheet.Cells.Item(5,4).NumberFormat = "@"
$newSheet.Cells.Item(5,4).Value2 = "3.91"
$newSheet.Cells.Item(5,4).ClearFormats()
$newSheet.Cells.Item(5,5).NumberFormat = "@"
$newSheet.Cells.Item(5,5).Value2 = "89.45"
$newSheet.Cells.Item(5,5).ClearFormats()
$newSheet.Cells.Item(5,6).NumberFormat = "@"
$newSheet.Cells.Item(5,6).Value2 = "5.13"
$newSheet.Cells.Item(5,6).ClearFormats()
$newSheet.Cells.Item(5,7).NumberFormat = "@"
$newSheet.Cells.Item(5,7).Value2 = "0.2006"
$newSheet.Cells.Item(5,7).ClearFormats()
$newSheet.Cells.Item(5,8).Value2 = 6

# row 6
$newSheet.Cells.Item(6,1).Value2 = 4
$newSheet.Cells.Item(6,2).NumberFormat = "@"
$newSheet.Cells.Item(6,2).Value2 = "862011"
$newSheet.Cells.Item(6,2).ClearFormats()
$newSheet.Cells.Item(6,3).NumberFormat = "@"
$newSheet.Cells.Item(6,3).Value2 = "光大阳光香港精选混合型集合资产管理计划（QDII）A 美元"
$newSheet.Cells.Item(6,3).ClearFormats()
$newSheet.Cells.Item(6,4).NumberFormat = "@"
$newSheet.Cells.Item(6,4).Value2 = "3.91"
$newSheet.Cells.Item(6,4).ClearFormats()
$newSheet.Cells.Item(6,5).NumberFormat = "@"
$newSheet.Cells.Item(6,5).Value2 = "89.45"
$newSheet.Cells.Item(6,5).ClearFormats()
$newSheet.Cells.Item(6,6).NumberFormat = "@"
$newSheet.Cells.Item(6,6).Value2 = "5.13"
$newSheet.Cells.Item(6,6).ClearFormats()
$newSheet.Cells.Item(6,7).NumberFormat = "@"
$newSheet.Cells.Item(6,7).Value2 = "0.2006"
$newSheet.Cells.Item(6,7).ClearFormats()
$newSheet.Cells.Item(6,8).Value2 = 6

# row 7
$newSheet.Cells.Item(7,1).Value2 = 5
$newSheet.Cells.Item(7,2).NumberFormat = "@"
$newSheet.Cells.Item(7,2).Value2 = "862012"
$newSheet.Cells.Item(7,2).ClearFormats()
$newSheet.Cells.Item(7,3).NumberFormat = "@"
$newSheet.Cells.Item(7,3).Value2 = "光大阳光香港精选混合型集合资产管理计划（QDII）C 人民币"
$newSheet.Cells.Item(7,3).ClearFormats()
$newSheet.Cells.Item(7,4).NumberFormat = "@"
$newSheet.Cells.Item(7,4).Value2 = "3.91"
$newSheet.Cells.Item(7,4).ClearFormats()
$newSheet.Cells.Item(7,5).NumberFormat = "@"
$newSheet.Cells.Item(7,5).Value2 = "89.45"
$newSheet.Cells.Item(7,5).ClearFormats()
$newSheet.Cells.Item(7,6).NumberFormat = "@"
$newSheet.Cells.Item(7,6).Value2 = "5.13"
$newSheet.Cells.Item(7,6).ClearFormats()
$newSheet.Cells.Item(7,7).NumberFormat = "@"
$newSheet.Cells.Item(7,7).Value2 = "0.2006"
$newSheet.Cells.Item(7,7).ClearFormats()
$newSheet.Cells.Item(7,8).Value2 = 6

# row 8
$newSheet.Cells.Item(8,1).Value2 = 6
$newSheet.Cells.Item(8,2).NumberFormat = "@"
$newSheet.Cells.Item(8,2).Value2 = "860018"
$newSheet.Cells.Item(8,2).ClearFormats()
$newSheet.Cells.Item(8,3).NumberFormat = "@"
$newSheet.Cells.Item(8,3).Value2 = "光大阳光智造混合A"
$newSheet.Cells.Item(8,3).ClearFormats()
$newSheet.Cells.Item(8,4).NumberFormat = "@"
$newSheet.Cells.Item(8,4).Value2 = "3.95"
$newSheet.Cells.Item(8,4).ClearFormats()
$newSheet.Cells.Item(8,5).NumberFormat = "@"
$newSheet.Cells.Item(8,5).Value2 = "90.37"
$newSheet.Cells.Item(8,5).ClearFormats()
$newSheet.Cells.Item(8,6).NumberFormat = "@"
$newSheet.Cells.Item(8,6).Value2 = "4.35"
$newSheet.Cells.Item(8,6).ClearFormats()
$newSheet.Cells.Item(8,7).NumberFormat = "@"
$newSheet.Cells.Item(8,7).Value2 = "0.1718"
$newSheet.Cells.Item(8,7).ClearFormats()
$newSheet.Cells.Item(8,8).Value2 = 8

# row 9
$newSheet.Cells.Item(9,1).Value2 = 7
$newSheet.Cells.Item(9,2).NumberFormat = "@"
$newSheet.Cells.Item(9,2).Value2 = "860038"
$newSheet.Cells.Item(9,2).ClearFormats()
$newSheet.Cells.Item(9,3).NumberFormat = "@"
$newSheet.Cells.Item(9,3).Value2 = "光大阳光智造混合B"
$newSheet.Cells.Item(9,3).ClearFormats()
$newSheet.Cells.Item(9,4).NumberFormat = "@"
$newSheet.Cells.Item(9,4).Value2 = "3.23"
$newSheet.Cells.Item(9,4).ClearFormats()
$newSheet.Cells.Item(9,5).NumberFormat = "@"
$newSheet.Cells.Item(9,5).Value2 = "90.37"
$newSheet.Cells.Item(9,5).ClearFormats()
$newSheet.Cells.Item(9,6).NumberFormat = "@"
$newSheet.Cells.Item(9,6).Value2 = "4.35"
$newSheet.Cells.Item(9,6).ClearFormats()
$newSheet.Cells.Item(9,7).NumberFormat = "@"
$newSheet.Cells.Item(9,7).Value2 = "0.1405"
$newSheet.Cells.Item(9,7).ClearFormats()
$newSheet.Cells.Item(9,8).Value2 = 8

# row 10
$newSheet.Cells.Item(10,1).Value2 = 8
$newSheet.Cells.Item(10,2).NumberFormat = "@"
$newSheet.Cells.Item(10,2).Value2 = "005646"
$newSheet.Cells.Item(10,2).ClearFormats()
$newSheet.Cells.Item(10,3).NumberFormat = "@"
$newSheet.Cells.Item(10,3).Value2 = "中海沪港深多策略灵活配置混合"
$newSheet.Cells.Item(10,3).ClearFormats()
$newSheet.Cells.Item(10,4).NumberFormat = "@"
$newSheet.Cells.Item(10,4).Value2 = "1.78"
$newSheet.Cells.Item(10,4).ClearFormats()
$newSheet.Cells.Item(10,5).NumberFormat = "@"
$newSheet.Cells.Item(10,5).Value2 = "88.15"
$newSheet.Cells.Item(10,5).ClearFormats()
$newSheet.Cells.Item(10,6).NumberFormat = "@"
$newSheet.Cells.Item(10,6).Value2 = "7.37"
$newSheet.Cells.Item(10,6).ClearFormats()
$newSheet.Cells.Item(10,7).NumberFormat = "@"
$newSheet.Cells.Item(10,7).Value2 = "0.1312"
$newSheet.Cells.Item(10,7).ClearFormats()
$newSheet.Cells.Item(10,8).Value2 = 1

# row 11
$newSheet.Cells.Item(11,1).Value2 = 9
$newSheet.Cells.Item(11,2).NumberFormat = "@"
$newSheet.Cells.Item(11,2).Value2 = "860007"
$newSheet.Cells.Item(11,2).ClearFormats()
$newSheet.Cells.Item(11,3).NumberFormat = "@"
$newSheet.Cells.Item(11,3).Value2 = "光大阳光价值30个月持有期混合A"
$newSheet.Cells.Item(11,3).ClearFormats()
$newSheet.Cells.Item(11,4).NumberFormat = "@"
$newSheet.Cells.Item(11,4).Value2 = "2.97"
$newSheet.Cells.Item(11,4).ClearFormats()
$newSheet.Cells.Item(11,5).NumberFormat = "@"
$newSheet.Cells.Item(11,5).Value2 = "90.90"
$newSheet.Cells.Item(11,5).ClearFormats()
$newSheet.Cells.Item(11,6).NumberFormat = "@"
$newSheet.Cells.Item(11,6).Value2 = "3.40"
$newSheet.Cells.Item(11,6).ClearFormats()
$newSheet.Cells.Item(11,7).NumberFormat = "@"
$newSheet.Cells.Item(11,7).Value2 = "0.1010"
$newSheet.Cells.Item(11,7).ClearFormats()
$newSheet.Cells.Item(11,8).Value2 = 9

# row 12
$newSheet.Cells.Item(12,1).Value2 = 10
$newSheet.Cells.Item(12,2).NumberFormat = "@"
$newSheet.Cells.Item(12,2).Value2 = "580008"
$newSheet.Cells.Item(12,2).ClearFormats()
$newSheet.Cells.Item(12,3).NumberFormat = "@"
$newSheet.Cells.Item(12,3).Value2 = "东吴新产业精选股票A"
$newSheet.Cells.Item(12,3).ClearFormats()
$newSheet.Cells.Item(12,4).NumberFormat = "@"
$newSheet.Cells.Item(12,4).Value2 = "2.32"
$newSheet.Cells.Item(12,4).ClearFormats()
$newSheet.Cells.Item(12,5).NumberFormat = "@"
$newSheet.Cells.Item(12,5).Value2 = "89.77"
$newSheet.Cells.Item(12,5).ClearFormats()
$newSheet.Cells.Item(12,6).NumberFormat = "@"
$newSheet.Cells.Item(12,6).Value2 = "3.87"
$newSheet.Cells.Item(12,6).ClearFormats()
$newSheet.Cells.Item(12,7).NumberFormat = "@"
$newSheet.Cells.Item(12,7).Value2 = "0.0898"
$newSheet.Cells.Item(12,7).ClearFormats()
$newSheet.Cells.Item(12,8).Value2 = 8

# row 13
$newSheet.Cells.Item(13,1).Value2 = 11
$newSheet.Cells.Item(13,2).NumberFormat = "@"
$newSheet.Cells.Item(13,2).Value2 = "012358"
$newSheet.Cells.Item(13,2).ClearFormats()
$newSheet.Cells.Item(13,3).NumberFormat = "@"
$newSheet.Cells.Item(13,3).Value2 = "汇丰晋信医疗先锋混合型证券投资基金A"
$newSheet.Cells.Item(13,3).ClearFormats()
$newSheet.Cells.Item(13,4).NumberFormat = "@"
$newSheet.Cells.Item(13,4).Value2 = "2.12"
$newSheet.Cells.Item(13,4).ClearFormats()
$newSheet.Cells.Item(13,5).NumberFormat = "@"
$newSheet.Cells.Item(13,5).Value2 = "85.40"
$newSheet.Cells.Item(13,5).ClearFormats()
$newSheet.Cells.Item(13,6).NumberFormat = "@"
$newSheet.Cells.Item(13,6).Value2 = "3.97"
$newSheet.Cells.Item(13,6).ClearFormats()
$newSheet.Cells.Item(13,7).NumberFormat = "@"
$newSheet.Cells.Item(13,7).Value2 = "0.0842"
$newSheet.Cells.Item(13,7).ClearFormats()
$newSheet.Cells.Item(13,8).Value2 = 6

# row 14
$newSheet.Cells.Item(14,1).Value2 = 12
$newSheet.Cells.Item(14,2).NumberFormat = "@"
$newSheet.Cells.Item(14,2).Value2 = "860027"
$newSheet.Cells.Item(14,2).ClearFormats()
$newSheet.Cells.Item(14,3).NumberFormat = "@"
$newSheet.Cells.Item(14,3).Value2 = "光大阳光价值30个月持有期混合B"
$newSheet.Cells.Item(14,3).ClearFormats()
$newSheet.Cells.Item(14,4).NumberFormat = "@"
$newSheet.Cells.Item(14,4).Value2 = "2.31"
$newSheet.Cells.Item(14,4).ClearFormats()
$newSheet.Cells.Item(14,5).NumberFormat = "@"
$newSheet.Cells.Item(14,5).Value2 = "90.90"
$newSheet.Cells.Item(14,5).ClearFormats()
$newSheet.Cells.Item(14,6).NumberFormat = "@"
$newSheet.Cells.Item(14,6).Value2 = "3.40"
$newSheet.Cells.Item(14,6).ClearFormats()
$newSheet.Cells.Item(14,7).NumberFormat = "@"
$newSheet.Cells.Item(14,7).Value2 = "0.0785"
$newSheet.Cells.Item(14,7).ClearFormats()
$newSheet.Cells.Item(14,8).Value2 = 9

# row 15
$newSheet.Cells.Item(15,1).Value2 = 13
$newSheet.Cells.Item(15,2).NumberFormat = "@"
$newSheet.Cells.Item(15,2).Value2 = "860008"
$newSheet.Cells.Item(15,2).ClearFormats()
$newSheet.Cells.Item(15,3).NumberFormat = "@"
$newSheet.Cells.Item(15,3).Value2 = "光大阳光生活 18 个月持有期混合型集合资产管理计划A"
$newSheet.Cells.Item(15,3).ClearFormats()
$newSheet.Cells.Item(15,4).NumberFormat = "@"
$newSheet.Cells.Item(15,4).Value2 = "0.51"
$newSheet.Cells.Item(15,4).ClearFormats()
$newSheet.Cells.Item(15,5).NumberFormat = "@"
$newSheet.Cells.Item(15,5).Value2 = "88.26"
$newSheet.Cells.Item(15,5).ClearFormats()
$newSheet.Cells.Item(15,6).NumberFormat = "@"
$newSheet.Cells.Item(15,6).Value2 = "4.34"
$newSheet.Cells.Item(15,6).ClearFormats()
$newSheet.Cells.Item(15,7).NumberFormat = "@"
$newSheet.Cells.Item(15,7).Value2 = "0.0221"
$newSheet.Cells.Item(15,7).ClearFormats()
$newSheet.Cells.Item(15,8).Value2 = 5

# row 16
$newSheet.Cells.Item(16,1).Value2 = 14
$newSheet.Cells.Item(16,2).NumberFormat = "@"
$newSheet.Cells.Item(16,2).Value2 = "860039"
$newSheet.Cells.Item(16,2).ClearFormats()
$newSheet.Cells.Item(16,3).NumberFormat = "@"
$newSheet.Cells.Item(16,3).Value2 = "光大阳光智造混合C"
$newSheet.Cells.Item(16,3).ClearFormats()
$newSheet.Cells.Item(16,4).NumberFormat = "@"
$newSheet.Cells.Item(16,4).Value2 = "0.26"
$newSheet.Cells.Item(16,4).ClearFormats()
$newSheet.Cells.Item(16,5).NumberFormat = "@"
$newSheet.Cells.Item(16,5).Value2 = "90.37"
$newSheet.Cells.Item(16,5).ClearFormats()
$newSheet.Cells.Item(16,6).NumberFormat = "@"
$newSheet.Cells.Item(16,6).Value2 = "4.35"
$newSheet.Cells.Item(16,6).ClearFormats()
$newSheet.Cells.Item(16,7).NumberFormat = "@"
$newSheet.Cells.Item(16,7).Value2 = "0.0113"
$newSheet.Cells.Item(16,7).ClearFormats()
$newSheet.Cells.Item(16,8).Value2 = 8

# row 17
$newSheet.Cells.Item(17,1).Value2 = 15
$newSheet.Cells.Item(17,2).NumberFormat = "@"
$newSheet.Cells.Item(17,2).Value2 = "860060"
$newSheet.Cells.Item(17,2).ClearFormats()
$newSheet.Cells.Item(17,3).NumberFormat = "@"
$newSheet.Cells.Item(17,3).Value2 = "光大阳光生活 18 个月持有期混合型集合资产管理计划B"
$newSheet.Cells.Item(17,3).ClearFormats()
$newSheet.Cells.Item(17,4).NumberFormat = "@"
$newSheet.Cells.Item(17,4).Value2 = "0.14"
$newSheet.Cells.Item(17,4).ClearFormats()
$newSheet.Cells.Item(17,5).NumberFormat = "@"
$newSheet.Cells.Item(17,5).Value2 = "88.26"
$newSheet.Cells.Item(17,5).ClearFormats()
$newSheet.Cells.Item(17,6).NumberFormat = "@"
$newSheet.Cells.Item(17,6).Value2 = "4.34"
$newSheet.Cells.Item(17,6).ClearFormats()
$newSheet.Cells.Item(17,7).NumberFormat = "@"
$newSheet.Cells.Item(17,7).Value2 = "0.0061"
$newSheet.Cells.Item(17,7).ClearFormats()
$newSheet.Cells.Item(17,8).Value2 = 5

# row 18
$newSheet.Cells.Item(18,1).Value2 = 16
$newSheet.Cells.Item(18,2).NumberFormat = "@"
$newSheet.Cells.Item(18,2).Value2 = "012359"
$newSheet.Cells.Item(18,2).ClearFormats()
$newSheet.Cells.Item(18,3).NumberFormat = "@"
$newSheet.Cells.Item(18,3).Value2 = "汇丰晋信医疗先锋混合型证券投资基金C"
$newSheet.Cells.Item(18,3).ClearFormats()
$newSheet.Cells.Item(18,4).NumberFormat = "@"
$newSheet.Cells.Item(18,4).Value2 = "0.11"
$newSheet.Cells.Item(18,4).ClearFormats()
$newSheet.Cells.Item(18,5).NumberFormat = "@"
$newSheet.Cells.Item(18,5).Value2 = "85.40"
$newSheet.Cells.Item(18,5).ClearFormats()
$newSheet.Cells.Item(18,6).NumberFormat = "@"
$newSheet.Cells.Item(18,6).Value2 = "3.97"
$newSheet.Cells.Item(18,6).ClearFormats()
$newSheet.Cells.Item(18,7).NumberFormat = "@"
$newSheet.Cells.Item(18,7).Value2 = "0.0044"
$newSheet.Cells.Item(18,7).ClearFormats()
$newSheet.Cells.Item(18,8).Value2 = 6

# row 19
$newSheet.Cells.Item(19,1).Value2 = 17
$newSheet.Cells.Item(19,2).NumberFormat = "@"
$newSheet.Cells.Item(19,2).Value2 = "860061"
$newSheet.Cells.Item(19,2).ClearFormats()
$newSheet.Cells.Item(19,3).NumberFormat = "@"
$newSheet.Cells.Item(19,3).Value2 = "光大阳光生活 18 个月持有期混合型集合资产管理计划C"
$newSheet.Cells.Item(19,3).ClearFormats()
$newSheet.Cells.Item(19,4).NumberFormat = "@"
$newSheet.Cells.Item(19,4).Value2 = "0.03"
$newSheet.Cells.Item(19,4).ClearFormats()
$newSheet.Cells.Item(19,5).NumberFormat = "@"
$newSheet.Cells.Item(19,5).Value2 = "88.26"
$newSheet.Cells.Item(19,5).ClearFormats()
$newSheet.Cells.Item(19,6).NumberFormat = "@"
$newSheet.Cells.Item(19,6).Value2 = "4.34"
$newSheet.Cells.Item(19,6).ClearFormats()
$newSheet.Cells.Item(19,7).NumberFormat = "@"
$newSheet.Cells.Item(19,7).Value2 = "0.0013"
$newSheet.Cells.Item(19,7).ClearFormats()
$newSheet.Cells.Item(19,8).Value2 = 5
